# Adds the next forecast date (2020-05-03 target-date column / 2020-05-17
# "created on" row) to both the "cases" and "deaths" sheets, matching a
# weekly-refresh upload of forecasts_table_RJ.xlsx.
#
# Layout recap: col A = "created on" date (string), row 1 = "target date"
# header (string, B1="Observed", C1.. = dates), data cells are plain
# numbers. New data lands in column X (24) and row 36.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, [string]$text) {
    # Writing a literal yyyy-mm-dd string via .Value auto-coerces to a
    # date serial (Excel's normal autocomplete behavior). The source
    # sheet stores these as plain text shared strings instead, so route
    # the literal through a text formula and flatten it back down to a
    # static value via copy/paste-values -- this keeps the cell's style
    # untouched (no NumberFormat/quotePrefix side effects).
    $cell = $ws.Cells.Item($row, $col)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$sheets = @(
    @{ Name = "cases";  B22 = 11139; X = @{23=11751;24=12559;25=13420;26=14395;27=15289;28=16267;29=17141;30=18114;31=19261;32=20249;33=21225;34=21778;35=22649;36=23222} },
    @{ Name = "deaths"; B22 = 1019;  X = @{23=1095; 24=1166; 25=1232; 26=1297; 27=1364; 28=1422; 29=1465; 30=1522; 31=1602; 32=1654; 33=1704; 34=1735; 35=1774; 36=1801 } }
)

foreach ($s in $sheets) {
    $ws = $wb.Worksheets.Item($s.Name)

    # New column header X1 = "2020-05-03" (next target date after W1).
    Set-TextCell $ws 1 24 "2020-05-03"

    # Fill in the new X column's observed/forecast values, rows 23-35
    # (row 22's diagonal start is blank, matching the existing staircase).
    foreach ($r in 23..35) {
        $ws.Cells.Item($r, 24).Value = $s.X[$r]
    }

    # B22 had been left blank; it now has its "observed" value filled in.
    $ws.Cells.Item(22, 2).Value = $s.B22

    # New row 36 = forecast created on 2020-05-17.
    Set-TextCell $ws 36 1 "2020-05-17"
    $ws.Cells.Item(36, 24).Value = $s.X[36]
}
